$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" metadata value (row 8, column B) to reflect the new
# generation timestamp.
$ws.Range("B8").Value = "2022-01-21T07:49:24+01:00"

# The IG now lists an additional author, so a second "Contact" row is
# added right after the existing one. Insert two new blank rows after
# the current Contact row (row 11) - this shifts every row below it
# down by two (rows 12-21 become rows 14-23).
$ws.Rows("12:13").Insert()

# Populate the two new rows with the same Contact / "No display for
# ContactDetail" values used by the existing Contact row, copying both
# the values and the cell formatting (borders, alignment, etc.) so the
# new rows are indistinguishable in style from the rest of the table.
$ws.Range("A10:B11").Copy()
$ws.Range("A12:B13").PasteSpecial(-4122)
$ws.Range("A10:B11").Copy()
$ws.Range("A12:B13").PasteSpecial(-4123)
